$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "294.24"
Set-TextValue $ws.Range("E2") "-6.25%"
Set-TextValue $ws.Range("D3") "40.67"
Set-TextValue $ws.Range("E3") "-0.49%"
Set-TextValue $ws.Range("D4") "5.023"
Set-TextValue $ws.Range("E4") "-2.19%"
Set-TextValue $ws.Range("D5") "0.07398"
Set-TextValue $ws.Range("E5") "-3.05%"
Set-TextValue $ws.Range("D6") "4.284"
Set-TextValue $ws.Range("E6") "-1.15%"
Set-TextValue $ws.Range("E7") "-8.52%"
Set-TextValue $ws.Range("D8") "0.9246"
Set-TextValue $ws.Range("E8") "-1.00%"
Set-TextValue $ws.Range("D9") "2.356"
Set-TextValue $ws.Range("E9") "-2.85%"
Set-TextValue $ws.Range("D10") "0.1146"
Set-TextValue $ws.Range("E10") "-8.72%"
Set-TextValue $ws.Range("D11") "0.1728"
Set-TextValue $ws.Range("E11") "-5.99%"
Set-TextValue $ws.Range("D12") "0.08685"
Set-TextValue $ws.Range("E12") "-3.93%"
Set-TextValue $ws.Range("E13") "0.25%"
Set-TextValue $ws.Range("D14") "0.1054"
Set-TextValue $ws.Range("E14") "-0.25%"
Set-TextValue $ws.Range("D15") "0.001265"
Set-TextValue $ws.Range("E15") "-1.58%"
Set-TextValue $ws.Range("D16") "0.005904"
Set-TextValue $ws.Range("E16") "0.76%"
Set-TextValue $ws.Range("D17") "3.415"
Set-TextValue $ws.Range("E17") "1.55%"
Set-TextValue $ws.Range("E18") "-2.22%"
Set-TextValue $ws.Range("D19") "7.681"
Set-TextValue $ws.Range("E19") "-8.90%"
Set-TextValue $ws.Range("D20") "0.1379"
Set-TextValue $ws.Range("E20") "2.32%"
Set-TextValue $ws.Range("E21") "5.10%"
Set-TextValue $ws.Range("D22") "0.03871"
Set-TextValue $ws.Range("E22") "-4.33%"
Set-TextValue $ws.Range("E23") "-0.64%"
Set-TextValue $ws.Range("D24") "0.003879"
Set-TextValue $ws.Range("E24") "-4.32%"
Set-TextValue $ws.Range("D25") "0.0001278"
Set-TextValue $ws.Range("E25") "0.23%"
Set-TextValue $ws.Range("D26") "0.0003718"
Set-TextValue $ws.Range("D38") "0.02337"
Set-TextValue $ws.Range("E38") "-5.74%"
Set-TextValue $ws.Range("D39") "0.05023"
Set-TextValue $ws.Range("E39") "-3.34%"
Set-TextValue $ws.Range("D40") "0.005989"
Set-TextValue $ws.Range("E40") "176.14%"
Set-TextValue $ws.Range("D41") "0.007674"
Set-TextValue $ws.Range("E41") "-1.54%"
Set-TextValue $ws.Range("D42") "0.1287"
Set-TextValue $ws.Range("E42") "-0.92%"
Set-TextValue $ws.Range("D43") "0.007339"
Set-TextValue $ws.Range("E43") "-0.23%"
Set-TextValue $ws.Range("D44") "0.007803"
Set-TextValue $ws.Range("E44") "-4.63%"
Set-TextValue $ws.Range("D45") "0.3163"
Set-TextValue $ws.Range("E45") "1.14%"
Set-TextValue $ws.Range("D46") "0.00006391"
Set-TextValue $ws.Range("E46") "-4.00%"
Set-TextValue $ws.Range("E47") "-0.47%"
Set-TextValue $ws.Range("D48") "0.01696"
Set-TextValue $ws.Range("E48") "-93.32%"
Set-TextValue $ws.Range("D49") "0.00002098"
Set-TextValue $ws.Range("E49") "-0.47%"
Set-TextValue $ws.Range("D50") "0.0001998"
Set-TextValue $ws.Range("E50") "-0.47%"
